$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 13: 50230552, thiago, 5, 4, 6, 7, 5.5, reprovado
Set-TextCell $ws.Range("A13") "50230552"
$ws.Range("B13").Value = "thiago"
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 6
$ws.Range("F13").Value = 7
$ws.Range("G13").Value = 5.5
$ws.Range("H13").Value = "reprovado"

# Row 14: 54023320, val, 2, 3, 42, 1, 12, aprovado
Set-TextCell $ws.Range("A14") "54023320"
$ws.Range("B14").Value = "val"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 42
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 12
$ws.Range("H14").Value = "aprovado"

# Row 15: 50230552, thiago, 5, 5, 5, 5, 5, reprovado
Set-TextCell $ws.Range("A15") "50230552"
$ws.Range("B15").Value = "thiago"
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = "reprovado"

# Row 16: 423, thiago, 5, 5, 5, 5, 5, reprovado
Set-TextCell $ws.Range("A16") "423"
$ws.Range("B16").Value = "thiago"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = "reprovado"

# Row 17: val, 5, 5, 5, 5, 5, 5, reprovado
$ws.Range("A17").Value = "val"
Set-TextCell $ws.Range("B17") "5"
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = "reprovado"
